# Generate Report for Handback
#
# Row 8 in both the "zh-cn" and "de-de" sheets represents the handback of
# 56e8fc5b-d7ec-48ac-aecb-2db5c3065f0d.md. A new handback came in, but it was
# not built on top of the latest handoff, so the report records the mismatch:
#   - "Latest Handback File" (I) / "Latest Target File" (J) / "Latest Handback
#     DateTime" (K) get filled in
#   - "Error Detail" (P) gets the "not the latest" message
#   - column P is widened so the message is readable
#   - a hyperlink is added on the new "Latest Handback File" cell

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cef6d20771a4880ecf28b857b1220d452565fc2d/e2e/56e8fc5b-d7ec-48ac-aecb-2db5c3065f0d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49caf348d4d2df4ee749b10c1dca4362ec2a3bf9/e2e/56e8fc5b-d7ec-48ac-aecb-2db5c3065f0d.md."
$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cef6d20771a4880ecf28b857b1220d452565fc2d/e2e/56e8fc5b-d7ec-48ac-aecb-2db5c3065f0d.md"
$handbackDisplay = "56e8fc5b-d7ec-48ac-aecb-2db5c3065f0d.md"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Columns.Item(16).ColumnWidth = 39.17

$wsZh.Range("I8").Value = $handbackDisplay
$wsZh.Range("J8").Value = "56e8fc5b-d7ec-48ac-aecb-2db5c3065f0d.7d4a89fe914e4a23568367f89586a96cbce15a50.zh-cn.xlf"
$wsZh.Range("K8").Value = "2016-08-18 20:44:34"
$wsZh.Range("P8").Value = $errorDetail

$wsZh.Hyperlinks.Add($wsZh.Range("I8"), $handbackUrl, "", "", $handbackDisplay)

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns.Item(16).ColumnWidth = 39.17

$wsDe.Range("I8").Value = $handbackDisplay
$wsDe.Range("J8").Value = "56e8fc5b-d7ec-48ac-aecb-2db5c3065f0d.7d4a89fe914e4a23568367f89586a96cbce15a50.de-de.xlf"
$wsDe.Range("K8").Value = "2016-08-18 20:44:43"
$wsDe.Range("P8").Value = $errorDetail

$wsDe.Hyperlinks.Add($wsDe.Range("I8"), $handbackUrl, "", "", $handbackDisplay)
